$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 210
$ws.Range("I11").Value = 210
$ws.Range("K11").Value = 210
$ws.Range("M11").Value = -70

$ws.Range("H17").Value = 9492.846
$ws.Range("J17").Value = 10158.917
$ws.Range("L17").Value = 30476.751
$ws.Range("N17").Value = -30812.751

$ws.Range("N37").ClearContents()
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0

$ws.Range("H86").Value = 6642.857
$ws.Range("I86").Value = 6750
$ws.Range("J86").Value = 6600
$ws.Range("K86").Value = 6750
$ws.Range("L86").Value = 6600
$ws.Range("M86").Value = -5627
$ws.Range("N86").Value = -8846

$ws.Range("H89").Value = 6642.857
$ws.Range("I89").Value = 6750
$ws.Range("J89").Value = 6600
$ws.Range("K89").Value = 33750
$ws.Range("L89").Value = 33000
$ws.Range("M89").Value = -28134
$ws.Range("N89").Value = -44232

$ws.Range("H98").Value = 187109.88
$ws.Range("I98").Value = 1100.3334
$ws.Range("K98").Value = 1100.3334
$ws.Range("M98").Value = 397.6666

$ws.Range("H122").Value = 187109.88
$ws.Range("I122").Value = 1100.3334
$ws.Range("K122").Value = 3301.0002
$ws.Range("M122").Value = -851.0001999999999

$ws.Range("H132").Value = 2225.6667
$ws.Range("J132").Value = 3999.5
$ws.Range("L132").Value = 11998.5
$ws.Range("N132").Value = -17058.5

$ws.Range("H141").Value = 2806.3
$ws.Range("I141").Value = 2806.3
$ws.Range("K141").Value = 8418.900000000001
$ws.Range("M141").Value = -3238.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -732

$ws.Range("H32").Value = 3454.5688
$ws.Range("I32").Value = 2266.0186
$ws.Range("K32").Value = 2266.0186
$ws.Range("M32").Value = -1979.0186

$ws.Range("H74").Value = 12348146
$ws.Range("I74").Value = 13890247
$ws.Range("K74").Value = 13890247
$ws.Range("M74").Value = -13889373

$ws.Range("H77").Value = 12348146
$ws.Range("I77").Value = 13890247
$ws.Range("K77").Value = 69451235
$ws.Range("M77").Value = -69446867

$ws.Range("H97").Value = 1204.1666
$ws.Range("I97").Value = 859.0909
$ws.Range("K97").Value = 859.0909
$ws.Range("M97").Value = -363.0909

$ws.Range("H119").Value = 67877.664
$ws.Range("J119").Value = 67877.664
$ws.Range("L119").Value = 67877.664
$ws.Range("N119").Value = -77553.664

$ws.Range("H122").Value = 4117.0527
$ws.Range("I122").Value = 3861.6365
$ws.Range("K122").Value = 11584.9095
$ws.Range("M122").Value = -9134.9095

$ws.Range("H125").Value = 212123
$ws.Range("J125").Value = 212123
$ws.Range("L125").Value = 212123
$ws.Range("N125").Value = -221963

$ws.Range("H132").Value = 2585.257
$ws.Range("I132").Value = 1968.7878
$ws.Range("K132").Value = 5906.3634
$ws.Range("M132").Value = -3376.3634

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2160.9167
$ws.Range("I94").Value = 2113.5625
$ws.Range("J94").Value = 2255.625
$ws.Range("K94").Value = 2113.5625
$ws.Range("L94").Value = 2255.625
$ws.Range("M94").Value = -1662.5625
$ws.Range("N94").Value = -3157.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23684.5
$ws.Range("I31").Value = 1491.8055
$ws.Range("K31").Value = 1491.8055
$ws.Range("M31").Value = -1196.8055

$ws.Range("H34").Value = 23684.5
$ws.Range("I34").Value = 1491.8055
$ws.Range("K34").Value = 1491.8055
$ws.Range("M34").Value = -1289.8055

$ws.Range("H109").Value = 70738.89
$ws.Range("J109").Value = 70738.89
$ws.Range("L109").Value = 70738.89
$ws.Range("N109").Value = -72818.89

$ws.Range("H122").Value = 4076
$ws.Range("J122").Value = 9312.875
$ws.Range("L122").Value = 27938.625
$ws.Range("N122").Value = -32838.625

$ws.Range("H134").Value = 2432.1482
$ws.Range("I134").Value = 1485.5834
$ws.Range("J134").Value = 10004.667
$ws.Range("K134").Value = 4456.7502
$ws.Range("L134").Value = 30014.001
$ws.Range("M134").Value = -1921.7502
$ws.Range("N134").Value = -35084.001

$ws.Range("H141").Value = 391336.84
$ws.Range("J141").Value = 391336.84
$ws.Range("L141").Value = 391336.84
$ws.Range("N141").Value = -401696.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 133422.08
$ws.Range("J2").Value = 145544.19
$ws.Range("L2").Value = 873265.14
$ws.Range("N2").Value = -873491.14

$ws.Range("H5").Value = 2851.8
$ws.Range("I5").Value = 416.33334
$ws.Range("J5").Value = 6505
$ws.Range("K5").Value = 1249.00002
$ws.Range("L5").Value = 19515
$ws.Range("M5").Value = -1137.00002
$ws.Range("N5").Value = -19739

$ws.Range("H8").Value = 266.5
$ws.Range("I8").Value = 266.5
$ws.Range("K8").Value = 799.5
$ws.Range("M8").Value = -660.5

$ws.Range("H70").Value = 14671.333
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685

$ws.Range("H73").Value = 14671.333
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908

$ws.Range("H116").Value = 3633.1428
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 3905.3333
$ws.Range("K116").Value = 6000
$ws.Range("L116").Value = 11715.9999
$ws.Range("M116").Value = -2558
$ws.Range("N116").Value = -18599.9999

$ws.Range("H122").Value = 2245.7368
$ws.Range("J122").Value = 3432.818
$ws.Range("L122").Value = 30895.362
$ws.Range("N122").Value = -35795.362

$ws.Range("H135").Value = 2851.8
$ws.Range("I135").Value = 416.33334
$ws.Range("J135").Value = 6505
$ws.Range("K135").Value = 3747.00006
$ws.Range("L135").Value = 58545
$ws.Range("M135").Value = -1212.00006
$ws.Range("N135").Value = -63615

$ws.Range("H140").Value = 2535.5
$ws.Range("I140").Value = 2535.5
$ws.Range("K140").Value = 7606.5
$ws.Range("M140").Value = -2426.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8654.729499999999
$ws.Range("I122").Value = 9821.237999999999
$ws.Range("J122").Value = 7123.6875
$ws.Range("K122").Value = 29463.714
$ws.Range("L122").Value = 21371.0625
$ws.Range("M122").Value = -27013.714
$ws.Range("N122").Value = -26271.0625

$ws.Range("H132").Value = 3038.7778
$ws.Range("I132").Value = 2258
$ws.Range("J132").Value = 3879.6155
$ws.Range("K132").Value = 6774
$ws.Range("L132").Value = 11638.8465
$ws.Range("M132").Value = -4244
$ws.Range("N132").Value = -16698.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2887.0356
$ws.Range("I22").Value = 1740.9445
$ws.Range("K22").Value = 1740.9445
$ws.Range("M22").Value = -1445.9445

$ws.Range("H27").Value = 2887.0356
$ws.Range("I27").Value = 1740.9445
$ws.Range("K27").Value = 1740.9445
$ws.Range("M27").Value = -1633.9445

$ws.Range("H40").Value = 7950.9653
$ws.Range("I40").Value = 6840.8945
$ws.Range("K40").Value = 6840.8945
$ws.Range("M40").Value = -6704.8945

$ws.Range("H122").Value = 190012.23
$ws.Range("I122").Value = 217497.16
$ws.Range("J122").Value = 15941
$ws.Range("K122").Value = 652491.48
$ws.Range("L122").Value = 47823
$ws.Range("M122").Value = -650041.48
$ws.Range("N122").Value = -52723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1684
$ws.Range("J4").Value = 1684
$ws.Range("L4").Value = 1684
$ws.Range("N4").Value = -1910

$ws.Range("H109").Value = 61710.332
$ws.Range("J109").Value = 61710.332
$ws.Range("L109").Value = 61710.332
$ws.Range("N109").Value = -64484.332

$ws.Range("H122").Value = 4900.1113
$ws.Range("I122").Value = 2399
$ws.Range("J122").Value = 6901
$ws.Range("K122").Value = 7197
$ws.Range("L122").Value = 20703
$ws.Range("M122").Value = -4747
$ws.Range("N122").Value = -25603

$ws.Range("H126").Value = 1305.0526
$ws.Range("I126").Value = 1083.4166
$ws.Range("J126").Value = 1685
$ws.Range("K126").Value = 3250.2498
$ws.Range("L126").Value = 5055
$ws.Range("M126").Value = -780.2498000000001
$ws.Range("N126").Value = -9995

$ws.Range("H132").Value = 3466.0667
$ws.Range("I132").Value = 1921.9231
$ws.Range("K132").Value = 5765.7693
$ws.Range("M132").Value = -3235.7693

$ws.Range("H136").Value = 3507.3667
$ws.Range("I136").Value = 2577.476
$ws.Range("J136").Value = 5677.1113
$ws.Range("K136").Value = 7732.428
$ws.Range("L136").Value = 17031.3339
$ws.Range("M136").Value = -5182.428
$ws.Range("N136").Value = -22131.3339
